$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New attribute rows (8-13), written column-by-column in the exact
# order the original author did (B/C pairs per row, except rows 10/11
# whose B cells were both entered before either C cell) so shared
# strings land at the same index the source workbook has.

# Row 8: 防御 (defense)
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "防御"
$ws.Range("C8").Value = "防御塔防御+{0}"
$ws.Range("D8").Value = "attr"
$ws.Range("E8").Value = 2

# Row 9: 魔力 (magic power)
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "魔力"
$ws.Range("C9").Value = "防御塔魔力+{0}"
$ws.Range("D9").Value = "attr"
$ws.Range("E9").Value = 2

# Rows 10 & 11: 命中 / 回避 (hit / dodge) - B cells first, then C cells
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "命中"
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "回避"
$ws.Range("C10").Value = "防御塔命中+{0}"
$ws.Range("C11").Value = "防御塔回避+{0}"
$ws.Range("D10").Value = "attr"
$ws.Range("E10").Value = 2
$ws.Range("D11").Value = "attr"
$ws.Range("E11").Value = 2

# Row 12: 暴击 (critical hit)
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "暴击"
$ws.Range("C12").Value = "防御塔暴击+{0}"
$ws.Range("D12").Value = "attr"
$ws.Range("E12").Value = 2

# Row 13: 幸运 (luck)
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "幸运"
$ws.Range("C13").Value = "防御塔幸运+{0}"
$ws.Range("D13").Value = "attr"
$ws.Range("E13").Value = 2

# --- Resize the table / autofilter to cover the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E13"))

# --- Column A now has an explicit width like columns B/C
$ws.Columns.Item(1).ColumnWidth = 8.2857142857143

# --- Move the selection to match where the author ended up (C13)
[void]$ws.Range("C13").Select()
